$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.023704012917921
$ws.Range("D2").Value = 1.049394706718483
$ws.Range("E2").Value = 1.024287885566793
$ws.Range("F2").Value = 1.05242664932019
$ws.Range("I2").Value = 1.037977565431155
$ws.Range("J2").Value = 1.028882740714751
$ws.Range("K2").Value = 1.052151885913943
$ws.Range("L2").Value = 1.027116827823794
$ws.Range("M2").Value = 1.055175413532968
$ws.Range("N2").Value = 1.013676984174211
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.024486955413413
$ws.Range("D3").Value = 1.050044013880634
$ws.Range("E3").Value = 1.024948177048208
$ws.Range("F3").Value = 1.053244022324577
$ws.Range("I3").Value = 1.038147354875485
$ws.Range("J3").Value = 1.02930513839509
$ws.Range("K3").Value = 1.052613851393742
$ws.Range("L3").Value = 1.027584578150941
$ws.Range("M3").Value = 1.055805618292446
$ws.Range("N3").Value = 1.01381815025647
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.024994320964799
$ws.Range("D4").Value = 1.05046463775939
$ws.Range("E4").Value = 1.025376467364823
$ws.Range("F4").Value = 1.053773920522564
$ws.Range("I4").Value = 1.038256342096704
$ws.Range("J4").Value = 1.029578539533939
$ws.Range("K4").Value = 1.052912564935757
$ws.Range("L4").Value = 1.027887598131265
$ws.Range("M4").Value = 1.056213752456323
$ws.Range("N4").Value = 1.013909484133034
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.025207795049034
$ws.Range("D5").Value = 1.050641580617352
$ws.Range("E5").Value = 1.025556767098993
$ws.Range("F5").Value = 1.053996927713056
$ws.Range("I5").Value = 1.038301949373271
$ws.Range("J5").Value = 1.029693495638608
$ws.Range("K5").Value = 1.053038092510945
$ws.Range("L5").Value = 1.028015071129478
$ws.Range("M5").Value = 1.056385413893836
$ws.Range("N5").Value = 1.013947878031158
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.025243648649691
$ws.Range("D6").Value = 1.050671296644462
$ws.Range("E6").Value = 1.025587054616055
$ws.Range("F6").Value = 1.0540343854939
$ws.Range("I6").Value = 1.038309594642772
$ws.Range("J6").Value = 1.029712798308662
$ws.Range("K6").Value = 1.053059166072092
$ws.Range("L6").Value = 1.028036479256366
$ws.Range("M6").Value = 1.05641424134969
$ws.Range("N6").Value = 1.013954324354138
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.024997172720342
$ws.Range("D7").Value = 1.050467001638711
$ws.Range("E7").Value = 1.025378875572585
$ws.Range("F7").Value = 1.053776899423169
$ws.Range("I7").Value = 1.038256952333231
$ws.Range("J7").Value = 1.029580075512854
$ws.Range("K7").Value = 1.052914242446062
$ws.Range("L7").Value = 1.027889301105775
$ws.Range("M7").Value = 1.056216045885616
$ws.Range("N7").Value = 1.013909997166016
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.023968455746908
$ws.Range("D8").Value = 1.049614042771164
$ws.Range("E8").Value = 1.024510818178264
$ws.Range("F8").Value = 1.052702675321351
$ws.Range("I8").Value = 1.03803512787527
$ws.Range("J8").Value = 1.0290254743314
$ws.Range("K8").Value = 1.052308051651017
$ws.Range("L8").Value = 1.027274832042936
$ws.Range("M8").Value = 1.055388320696573
$ws.Range("N8").Value = 1.013724693707391
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.022161544357699
$ws.Range("D9").Value = 1.048114776451254
$ws.Range("E9").Value = 1.022989221290179
$ws.Range("F9").Value = 1.05081753512059
$ws.Range("I9").Value = 1.037637555695681
$ws.Range("J9").Value = 1.028048877206913
$ws.Range("K9").Value = 1.051238334209697
$ws.Range("L9").Value = 1.026194834354435
$ws.Range("M9").Value = 1.053932521941787
$ws.Range("N9").Value = 1.013398109177758
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.020960957163914
$ws.Range("D10").Value = 1.04711791235629
$ws.Range("E10").Value = 1.021980337469564
$ws.Range("F10").Value = 1.049566136307951
$ws.Range("I10").Value = 1.037368056887984
$ws.Range("J10").Value = 1.027398346157633
$ws.Range("K10").Value = 1.050524252634974
$ws.Range("L10").Value = 1.0254767860287
$ws.Range("M10").Value = 1.052963954750088
$ws.Range("N10").Value = 1.013180376197445
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.020442066030526
$ws.Range("D11").Value = 1.04668691270539
$ws.Range("E11").Value = 1.021544811787307
$ws.Range("F11").Value = 1.049025566107074
$ws.Range("I11").Value = 1.037250315167748
$ws.Range("J11").Value = 1.02711680167035
$ws.Range("K11").Value = 1.050214845037507
$ws.Range("L11").Value = 1.025166343560417
$ws.Range("M11").Value = 1.052545044010489
$ws.Range("N11").Value = 1.013086098935921
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.020249474492886
$ws.Range("D12").Value = 1.04652691985417
$ws.Range("E12").Value = 1.021383239608566
$ws.Range("F12").Value = 1.048824970887943
$ws.Range("I12").Value = 1.037206424144171
$ws.Range("J12").Value = 1.027012245807087
$ws.Range("K12").Value = 1.050099887993356
$ws.Range("L12").Value = 1.025051104548288
$ws.Range("M12").Value = 1.052389516711111
$ws.Range("N12").Value = 1.013051081023612
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.020290779316611
$ws.Range("D13").Value = 1.046561234322763
$ws.Range("E13").Value = 1.021417888241453
$ws.Range("F13").Value = 1.048867990309658
$ws.Range("I13").Value = 1.0372158459873
$ws.Range("J13").Value = 1.027034672364925
$ws.Range("K13").Value = 1.050124547969171
$ws.Range("L13").Value = 1.025075820374439
$ws.Range("M13").Value = 1.052422874429385
$ws.Range("N13").Value = 1.0130585924381
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.020426143329562
$ws.Range("D14").Value = 1.046673685602307
$ws.Range("E14").Value = 1.021531452059601
$ws.Range("F14").Value = 1.049008980819627
$ws.Range("I14").Value = 1.037246690310168
$ws.Range("J14").Value = 1.027108158587752
$ws.Range("K14").Value = 1.050205343240567
$ws.Range("L14").Value = 1.025156816367251
$ws.Range("M14").Value = 1.052532186538504
$ws.Range("N14").Value = 1.013083204324063
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02050956519678
$ws.Range("D15").Value = 1.046742983808741
$ws.Range("E15").Value = 1.021601449229596
$ws.Range("F15").Value = 1.04909587584781
$ws.Range("I15").Value = 1.037265673800462
$ws.Range("J15").Value = 1.027153438868155
$ws.Range("K15").Value = 1.050255120040936
$ws.Range("L15").Value = 1.025206730407131
$ws.Range("M15").Value = 1.052599547305363
$ws.Range("N15").Value = 1.013098368647538
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.020995414875685
$ws.Range("D16").Value = 1.047146530237703
$ws.Range("E16").Value = 1.022009270016637
$ws.Range("F16").Value = 1.049602039621798
$ws.Range("I16").Value = 1.037375849027154
$ws.Range("J16").Value = 1.027417034393291
$ws.Range("K16").Value = 1.050544782803601
$ws.Range("L16").Value = 1.02549739925065
$ws.Range("M16").Value = 1.05299176686188
$ws.Range("N16").Value = 1.013186633163317
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.021300436820693
$ws.Range("D17").Value = 1.047399839506207
$ws.Range("E17").Value = 1.022265442078388
$ws.Range("F17").Value = 1.049919890746407
$ws.Range("I17").Value = 1.037444679291006
$ws.Range("J17").Value = 1.027582419234415
$ws.Range("K17").Value = 1.050726426694587
$ws.Range("L17").Value = 1.02567985691948
$ws.Range("M17").Value = 1.053237927081987
$ws.Range("N17").Value = 1.013242000182828
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.021478444528505
$ws.Range("D18").Value = 1.047547653003548
$ws.Range("E18").Value = 1.022414990875913
$ws.Range("F18").Value = 1.050105412606226
$ws.Range("I18").Value = 1.037484725750232
$ws.Range("J18").Value = 1.027678898810883
$ws.Range("K18").Value = 1.050832356562037
$ws.Range("L18").Value = 1.025786327336031
$ws.Range("M18").Value = 1.053381554865582
$ws.Range("N18").Value = 1.013274295043252
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.021539156371388
$ws.Range("D19").Value = 1.047598064121428
$ws.Range("E19").Value = 1.022466004796363
$ws.Range("F19").Value = 1.05016869184487
$ws.Range("I19").Value = 1.037498363390919
$ws.Range("J19").Value = 1.027711798085988
$ws.Range("K19").Value = 1.050868472478418
$ws.Range("L19").Value = 1.025822638745615
$ws.Range("M19").Value = 1.053430536090101
$ws.Range("N19").Value = 1.013285306777997
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.021267701161233
$ws.Range("D20").Value = 1.047372655340613
$ws.Range("E20").Value = 1.022237943989914
$ws.Range("F20").Value = 1.04988577541674
$ws.Range("I20").Value = 1.037437304903482
$ws.Range("J20").Value = 1.027564673611094
$ws.Range("K20").Value = 1.050706940075216
$ws.Range("L20").Value = 1.025660276175588
$ws.Range("M20").Value = 1.053211511576753
$ws.Range("N20").Value = 1.013236059799982
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02038627790863
$ws.Range("D21").Value = 1.046640568736758
$ws.Range("E21").Value = 1.021498004768697
$ws.Range("F21").Value = 1.048967457182715
$ws.Range("I21").Value = 1.037237611739865
$ws.Range("J21").Value = 1.02708651810084
$ws.Range("K21").Value = 1.050181551842523
$ws.Range("L21").Value = 1.025132963035398
$ws.Range("M21").Value = 1.052499994766774
$ws.Range("N21").Value = 1.01307595670951
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.019832947778163
$ws.Range("D22").Value = 1.046180854175439
$ws.Range("E22").Value = 1.021033942127719
$ws.Range("F22").Value = 1.048391212771472
$ws.Range("I22").Value = 1.037111151616179
$ws.Range("J22").Value = 1.026786012388898
$ws.Range("K22").Value = 1.049851050758103
$ws.Range("L22").Value = 1.024801844234825
$ws.Range("M22").Value = 1.052053069624916
$ws.Range("N22").Value = 1.012975298754973
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.020126196827173
$ws.Range("D23").Value = 1.046424502139066
$ws.Range("E23").Value = 1.021279839280239
$ws.Range("F23").Value = 1.048696582148432
$ws.Range("I23").Value = 1.037178276050709
$ws.Range("J23").Value = 1.026945303430898
$ws.Range("K23").Value = 1.050026271089048
$ws.Range("L23").Value = 1.024977335911459
$ws.Range("M23").Value = 1.052289951353618
$ws.Range("N23").Value = 1.013028658788145
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.021282492707784
$ws.Range("D24").Value = 1.047384938502123
$ws.Range("E24").Value = 1.022250368797324
$ws.Range("F24").Value = 1.049901190280809
$ws.Range("I24").Value = 1.037440637384166
$ws.Range("J24").Value = 1.027572692053199
$ws.Range("K24").Value = 1.050715745301542
$ws.Range("L24").Value = 1.025669123729141
$ws.Range("M24").Value = 1.053223447462506
$ws.Range("N24").Value = 1.013238744002689
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.022627973542318
$ws.Range("D25").Value = 1.048501915093876
$ws.Range("E25").Value = 1.023381627025873
$ws.Range("F25").Value = 1.051303953635202
$ws.Range("I25").Value = 1.037741125524305
$ws.Range("J25").Value = 1.028301262454842
$ws.Range("K25").Value = 1.051515053722637
$ws.Range("L25").Value = 1.026473702078231
$ws.Range("M25").Value = 1.054308542501598
$ws.Range("N25").Value = 1.013482542896089
